$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

$data = @(
    ,@("2025-07-25T09:57:35+00:00", "EXTERNAL:- GREAT! TV Week 30-33 Schedules", "schedules@globallistings-distribution.com", $true)
    ,@("2025-07-25T09:57:04+00:00", "EXTERNAL:- GREAT! TV Week 33 Schedule", "schedules@globallistings-distribution.com", $true)
    ,@("2025-07-25T09:56:57+00:00", "EXTERNAL:- ITV4 Post Press Change - Wednesday 30th July", "itv2-4-scheduling@itv.com", $false)
    ,@("2025-07-25T09:56:04+00:00", "EXTERNAL:- Cúla4 Schedule Week 34 English", "Sceideal@europus.ie", $true)
    ,@("2025-07-25T09:55:42+00:00", "EXTERNAL:- TNT Sports - Updated Schedules - Weeks 30 & 31", "James.Prestridge@wbd.com", $true)
    ,@("2025-07-25T09:54:22+00:00", "EXTERNAL:- Sceideal Cúla4 Seachtain 34 Gaeilge Béarla ", "Sceideal@europus.ie", $true)
    ,@("2025-07-25T09:35:47+00:00", "EXTERNAL:- Quest Red programme schedule Week 34", "ellie.tangokpho@wbd.com", $true)
    ,@("2025-07-25T09:35:45+00:00", "EXTERNAL:- Die-hard Football Fans Roll Up as RTÉ’s Up for the Match Returns Before Kerry and Donegal Take Centre Stage at Croke Park", "releases2@rte.ie", $false)
    ,@("2025-07-25T09:34:32+00:00", "U&W - Week 33 Listings", "Manisha.Anand@pa.media", $true)
    ,@("2025-07-25T09:34:28+00:00", "U&Alibi - Week 33 Listings", "Manisha.Anand@pa.media", $true)
    ,@("2025-07-25T09:34:05+00:00", "EXTERNAL:- Schedule for channel: LBC London from 28/07/2025 to 04/08/2025", "epgcaster@myradio.thisisglobal.com", $true)
    ,@("2025-07-25T09:31:11+00:00", "EXTERNAL:- BBC Radio London - Wk33 - 2025-08-22 - Friday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:31:08+00:00", "EXTERNAL:- BBC Radio London - Wk33 - 2025-08-21 - Thursday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:31:05+00:00", "EXTERNAL:- BBC Radio London - Wk33 - 2025-08-20 - Wednesday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:31:03+00:00", "EXTERNAL:- BBC Radio London - Wk33 - 2025-08-19 - Tuesday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:30:58+00:00", "EXTERNAL:- BBC Radio London - Wk33 - 2025-08-18 - Monday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:30:55+00:00", "EXTERNAL:- BBC Radio London - Wk33 - 2025-08-16 - Saturday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:30:48+00:00", "EXTERNAL:- BBC Radio London - Wk32 - 2025-08-15 - Friday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:30:47+00:00", "EXTERNAL:- BBC Radio London - Wk32 - 2025-08-13 - Wednesday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:30:47+00:00", "EXTERNAL:- BBC Radio London - Wk31 - 2025-08-08 - Friday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:30:46+00:00", "EXTERNAL:- BBC Radio London - Wk32 - 2025-08-14 - Thursday", "pressportal@bbc.co.uk", $true)
    ,@("2025-07-25T09:30:46+00:00", "EXTERNAL:- BBC Radio London - Wk32 - 2025-08-09 - Saturday", "pressportal@bbc.co.uk", $true)
)

foreach ($row in $data) {
    $newRow = $tbl.ListRows.Add()
    $r = $newRow.Range.Row
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
